$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$primer = $ws.Range("Z1")
$primer.WrapText = $true
$primer.WrapText = $false
$primer.Clear() | Out-Null

$ws.Range("D2").Value = "https://ecoinfaeet.github.io/new_website/"

$ws.Range("E1").Value = "titulo"

$ws.Range("D9").Value = "https://doi.org/10.7818/ECOS.1591"

# Plain (no wrap) titulo cells
$ws.Range("E11").Value = "Cómo escribir funciones en R"
$ws.Range("E12").Value = "Cómo crear paquetes en R"
$ws.Range("E14").Value = "Estadística circular aplicada a la Ecología"

# Wrapped titulo cells (style index 2: wrapText)
$ws.Range("E4").Value = "rOpenSci: cómo acceder `nde forma reproducible a repositorios `nde datos públicos"
$ws.Range("E4").WrapText = $true
$ws.Range("E5").Value = "Ajuste, interpretación y presentación de `nmodelos lineales: el valor p no es suficiente"
$ws.Range("E5").WrapText = $true
$ws.Range("E6").Value = "Procesadores de texto Markup: `nmás allá de MS Word"
$ws.Range("E6").WrapText = $true
$ws.Range("E7").Value = "Sobre el intercambio de datos de vegetación: `nel estándar ‘Veg-X’ y el paquete de R ‘VegX’"
$ws.Range("E7").WrapText = $true
$ws.Range("E8").Value = "¿Por qué usar GitHub? Diez pasos para `ndisfrutar de GitHub y no morir en el intento"
$ws.Range("E8").WrapText = $true
$ws.Range("E9").Value = "Ventajas de la estadística bayesiana frente `na la frecuentista: ¿por qué nos resistimos a usarla?"
$ws.Range("E9").WrapText = $true
$ws.Range("E10").Value = "Inferencia estadística a partir de varios `nmodelos y su utilidad en ecología"
$ws.Range("E10").WrapText = $true
$ws.Range("E13").Value = "Compartiendo datos en Ecología: `ncómo añadir más valor a los datos"
$ws.Range("E13").WrapText = $true
$ws.Range("E15").Value = "Quince consejos para mejorar nuestro `ncódigo y flujo de trabajo con R"
$ws.Range("E15").WrapText = $true
$ws.Range("E16").Value = "Cómo aplicar la cienciometría `na la investigación ecológica"
$ws.Range("E16").WrapText = $true
$ws.Range("E17").Value = "¡Se puede entender cómo `nfuncionan Git y GitHub!"
$ws.Range("E17").WrapText = $true
$ws.Range("E18").Value = "La unión hace la fuerza: `nmodelos de distribución de especies `nintegrando diferentes fuentes de datos"
$ws.Range("E18").WrapText = $true
$ws.Range("E19").Value = "Camelot: Una herramienta intuitiva `npara el manejo y procesamiento de imágenes de `ncámaras trampa utilizando inteligencia artificial"
$ws.Range("E19").WrapText = $true
$ws.Range("E20").Value = "Tidyverse: colección de paquetes `nde R para la ciencia de datos"
$ws.Range("E20").WrapText = $true
$ws.Range("E21").Value = "Búsqueda, descarga y limpieza `nde datos desde GBIF"
$ws.Range("E21").WrapText = $true
$ws.Range("E22").Value = "Estadística bayesiana"
$ws.Range("E22").WrapText = $true
$ws.Range("E23").Value = "Análisis de la ciencia ciudadana `nmediante modelos de ocupación"
$ws.Range("E23").WrapText = $true
$ws.Range("E24").Value = "Generación de visores de `ndatos espaciales con R"
$ws.Range("E24").WrapText = $true
$ws.Range("E25").Value = "Análisis de dinámicas poblacionales en R"
$ws.Range("E25").WrapText = $true
$ws.Range("E26").Value = "Introducción a Python"
$ws.Range("E26").WrapText = $true
$ws.Range("E27").Value = "Integración de datos en la estima `nde la distribución y abundancia animal"
$ws.Range("E27").WrapText = $true
$ws.Range("E28").Value = "Introducción al análisis `nespacial de patrones de puntos"
$ws.Range("E28").WrapText = $true
$ws.Range("E29").Value = "Introducción al mundo de la bioacústica"
$ws.Range("E29").WrapText = $true
$ws.Range("E30").Value = "Introducción a Zotero"
$ws.Range("E30").WrapText = $true
$ws.Range("E31").Value = "Introducción al uso de filogenias"
$ws.Range("E31").WrapText = $true
$ws.Range("E32").Value = "Regímenes dinámicos ecológicos"
$ws.Range("E32").WrapText = $true
$ws.Range("E33").Value = "labeleR: genera tus certificados y etiquetas"
$ws.Range("E33").WrapText = $true
$ws.Range("E34").Value = "¿Qué información puedo obtener `nde los datos PNOA-LiDAR?"
$ws.Range("E34").WrapText = $true
$ws.Range("E35").Value = "Inferir rango de distribución a `npartir de diversidad genética"
$ws.Range("E35").WrapText = $true
$ws.Range("E36").Value = "El papel de la IA en la ecología, `ncómo transformar imágenes en datos"
$ws.Range("E36").WrapText = $true

# Apply the existing Hyperlink cell style to D2:D23
$ws.Range("D2:D23").Style = "Hipervínculo"

# Hyperlinks, added in the exact order to reproduce rId2..rId21
$ws.Hyperlinks.Add($ws.Range("D2"), "https://ecoinfaeet.github.io/new_website/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://ecoinf.quarto.pub/iecoinf/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://doi.org/10.7818/ECOS.2017.26-2.08") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "https://doi.org/10.7818/ECOS.1570") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "https://doi.org/10.7818/ECOS.2017.26-3.14") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "https://doi.org/10.7818/ECOS.1604") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "https://doi.org/10.7818/ECOS.1591") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "https://doi.org/10.7818/ECOS.1699") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), "https://doi.org/10.7818/ECOS.1948") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D14"), "https://doi.org/10.7818/ECOS.1995") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D15"), "https://doi.org/10.7818/ECOS.2129") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D17"), "https://doi.org/10.7818/ECOS.2332") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D19"), "https://doi.org/10.7818/ECOS.2797") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D21"), "https://www.youtube.com/watch?v=VSUEi50tkAI") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D22"), "https://www.youtube.com/watch?v=usB7reMJxLU") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "https://doi.org/10.7818/ECOS.1880") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D13"), "https://doi.org/10.7818/ECOS.1838") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D16"), "https://doi.org/10.7818/ECOS.2256") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D18"), "https://doi.org/10.7818/ECOS.2527") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D23"), "https://www.youtube.com/watch?v=OKEtldANpHA") | Out-Null

$ws.Range("D2:D23").Style = "Hipervínculo"

# Row heights
$ws.Rows.Item(4).RowHeight = 145
$ws.Rows.Item(5).RowHeight = 159.5
$ws.Rows.Item(6).RowHeight = 101.5
$ws.Rows.Item(7).RowHeight = 188.5
$ws.Rows.Item(8).RowHeight = 174
$ws.Rows.Item(9).RowHeight = 188.5
$ws.Rows.Item(10).RowHeight = 145
$ws.Rows.Item(13).RowHeight = 130.5
$ws.Rows.Item(15).RowHeight = 130.5
$ws.Rows.Item(16).RowHeight = 116
$ws.Rows.Item(17).RowHeight = 101.5
$ws.Rows.Item(18).RowHeight = 188.5
$ws.Rows.Item(19).RowHeight = 275.5
$ws.Rows.Item(20).RowHeight = 116
$ws.Rows.Item(21).RowHeight = 116
$ws.Rows.Item(22).RowHeight = 58
$ws.Rows.Item(23).RowHeight = 145
$ws.Rows.Item(24).RowHeight = 101.5
$ws.Rows.Item(25).RowHeight = 101.5
$ws.Rows.Item(26).RowHeight = 43.5
$ws.Rows.Item(27).RowHeight = 145
$ws.Rows.Item(28).RowHeight = 116
$ws.Rows.Item(29).RowHeight = 87
$ws.Rows.Item(30).RowHeight = 43.5
$ws.Rows.Item(31).RowHeight = 58
$ws.Rows.Item(32).RowHeight = 87
$ws.Rows.Item(33).RowHeight = 87
$ws.Rows.Item(34).RowHeight = 130.5
$ws.Rows.Item(35).RowHeight = 116
$ws.Rows.Item(36).RowHeight = 130.5

# Column widths
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Final selection / view
$ws.Range("E37").Select() | Out-Null

